$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.500.07"
$ws.Range("E2").Value = "  +1.57%  "
$ws.Range("D3").Value = "1.571.58"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  -1.48%  "
$ws.Range("D5").Value = "'211.29"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("E7").Value = "  -1.57%  "
$ws.Range("D8").Value = "'22.92"
$ws.Range("E8").Value = "  +3.98%  "
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("E10").Value = "  -0.26%  "
$ws.Range("E11").Value = "  +1.28%  "
$ws.Range("D12").Value = "1.796.77"
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("D13").Value = "1.570.50"
$ws.Range("E13").Value = "  +0.38%  "
$ws.Range("D14").Value = "'3.76"
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("D16").Value = "27.467.79"
$ws.Range("E16").Value = "  +1.49%  "
$ws.Range("D17").Value = "'62.37"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").Value = "'226.18"
$ws.Range("E18").Value = "  +4.97%  "
$ws.Range("D19").Value = "'7.51"
$ws.Range("E19").Value = "  +1.33%  "
$ws.Range("E21").Value = "  -1.54%  "
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("E23").Value = "  +2.42%  "
$ws.Range("D24").Value = "'1.96"
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("D25").Value = "'150.51"
$ws.Range("E25").Value = "  -2.31%  "
$ws.Range("D26").Value = "'15.18"
$ws.Range("E26").Value = "  +0.88%  "
$ws.Range("E28").Value = "  +1.47%  "
$ws.Range("D29").Value = "'0.993"
$ws.Range("E29").Value = "  -1.45%  "
$ws.Range("D31").Value = "'0.0473"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").Value = "'3.25"
$ws.Range("E32").Value = "  +0.48%  "
$ws.Range("D33").Value = "1.456.17"
$ws.Range("E33").Value = "  +1.95%  "
$ws.Range("E34").Value = "  -1.92%  "
$ws.Range("E35").Value = "  +3.19%  "
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("E37").Value = "  -1.26%  "
$ws.Range("E38").Value = "  +0.49%  "
$ws.Range("D39").Value = "'0.541"
$ws.Range("E39").Value = "  +1.87%  "
$ws.Range("D40").Value = "'0.814"
$ws.Range("E40").Value = "  +0.20%  "
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "'0.993"
$ws.Range("E42").Value = "  -1.52%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.66"
$ws.Range("E43").Value = "  -3.12%  "
$ws.Range("D44").Value = "'1.85"
$ws.Range("E44").Value = "  +6.71%  "
$ws.Range("D45").Value = "'0.972"
$ws.Range("E45").Value = "  -2.97%  "
$ws.Range("D46").Value = "'63.93"
$ws.Range("E46").Value = "  -1.11%  "
$ws.Range("D47").Value = "1.708.01"
$ws.Range("E47").Value = "  +0.28%  "
$ws.Range("D48").Value = "'86.92"
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("E49").Value = "  +1.36%  "
$ws.Range("D50").Value = "'0.0526"
$ws.Range("E50").Value = "  +1.65%  "
$ws.Range("D51").Value = "'0.0948"
$ws.Range("E51").Value = "  -1.62%  "
